$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, extend column A styling (copy format from A7, which already has the bordered/bold/centered style)
# down through A25, so newly added rows 8-25 match the style of existing column-A header cells.
$ws.Cells.Item(7,1).Copy() | Out-Null
$ws.Range("A8:A25").PasteSpecial(-4122) | Out-Null

# Populate data rows 2-25 with the new timeseries values (columns A-E)
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = 6.35833331
$ws.Cells.Item(2,3).Value = 50.734
$ws.Cells.Item(2,4).Value = 6.35833331
$ws.Cells.Item(2,5).Value = 35.5138

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = 5.54166669
$ws.Cells.Item(3,3).Value = 48.092
$ws.Cells.Item(3,4).Value = 5.54166669
$ws.Cells.Item(3,5).Value = 33.6644

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = 5.54166669
$ws.Cells.Item(4,3).Value = 60.06700000000001
$ws.Cells.Item(4,4).Value = 5.54166669
$ws.Cells.Item(4,5).Value = 42.0469

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = 5.54166669
$ws.Cells.Item(5,3).Value = 59.18
$ws.Cells.Item(5,4).Value = 5.54166669
$ws.Cells.Item(5,5).Value = 41.426

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = 19.075
$ws.Cells.Item(6,3).Value = 59.261
$ws.Cells.Item(6,4).Value = 19.075
$ws.Cells.Item(6,5).Value = 41.4827

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = 29.16666669
$ws.Cells.Item(7,3).Value = 30.006
$ws.Cells.Item(7,4).Value = 29.16666669
$ws.Cells.Item(7,5).Value = 21.0042

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = 32.08333331
$ws.Cells.Item(8,3).Value = 33.838
$ws.Cells.Item(8,4).Value = 32.08333331
$ws.Cells.Item(8,5).Value = 23.6866

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = 29.16666669
$ws.Cells.Item(9,3).Value = 33.446
$ws.Cells.Item(9,4).Value = 29.16666669
$ws.Cells.Item(9,5).Value = 23.4122

$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = 24.5
$ws.Cells.Item(10,3).Value = 54.44
$ws.Cells.Item(10,4).Value = 24.5
$ws.Cells.Item(10,5).Value = 38.108

$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = 25.08333331
$ws.Cells.Item(11,3).Value = 58.648
$ws.Cells.Item(11,4).Value = 25.08333331
$ws.Cells.Item(11,5).Value = 41.0536

$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = 28.875
$ws.Cells.Item(12,3).Value = 41.757
$ws.Cells.Item(12,4).Value = 28.875
$ws.Cells.Item(12,5).Value = 29.2299

$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = 31.09166669
$ws.Cells.Item(13,3).Value = 26.659
$ws.Cells.Item(13,4).Value = 31.09166669
$ws.Cells.Item(13,5).Value = 18.6613

$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = 40.30833330999999
$ws.Cells.Item(14,3).Value = 26.996
$ws.Cells.Item(14,4).Value = 40.30833330999999
$ws.Cells.Item(14,5).Value = 18.8972

$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = 30.275
$ws.Cells.Item(15,3).Value = 25.118
$ws.Cells.Item(15,4).Value = 30.275
$ws.Cells.Item(15,5).Value = 17.5826

$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = 24.38333331
$ws.Cells.Item(16,3).Value = 18.521
$ws.Cells.Item(16,4).Value = 24.38333331
$ws.Cells.Item(16,5).Value = 12.9647

$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = 23.15833331
$ws.Cells.Item(17,3).Value = 15.541
$ws.Cells.Item(17,4).Value = 23.15833331
$ws.Cells.Item(17,5).Value = 10.8787

$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = 23.85833331
$ws.Cells.Item(18,3).Value = 51.249
$ws.Cells.Item(18,4).Value = 23.85833331
$ws.Cells.Item(18,5).Value = 35.8743

$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = 38.38333331
$ws.Cells.Item(19,3).Value = 88.39700000000001
$ws.Cells.Item(19,4).Value = 38.38333331
$ws.Cells.Item(19,5).Value = 61.8779

$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = 70
$ws.Cells.Item(20,3).Value = 86.38800000000001
$ws.Cells.Item(20,4).Value = 70
$ws.Cells.Item(20,5).Value = 60.4716

$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = 58.50833331
$ws.Cells.Item(21,3).Value = 88.33800000000001
$ws.Cells.Item(21,4).Value = 58.50833331
$ws.Cells.Item(21,5).Value = 61.8366

$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = 39.43333331
$ws.Cells.Item(22,3).Value = 95.18799999999999
$ws.Cells.Item(22,4).Value = 39.43333331
$ws.Cells.Item(22,5).Value = 66.63159999999999

$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = 28
$ws.Cells.Item(23,3).Value = 94.732
$ws.Cells.Item(23,4).Value = 28
$ws.Cells.Item(23,5).Value = 66.3124

$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = 17.5
$ws.Cells.Item(24,3).Value = 97.48400000000001
$ws.Cells.Item(24,4).Value = 17.5
$ws.Cells.Item(24,5).Value = 68.2388

$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = 11.9
$ws.Cells.Item(25,3).Value = 94
$ws.Cells.Item(25,4).Value = 11.9
$ws.Cells.Item(25,5).Value = 65.8

